$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin-list values (prices, 1h volume %, and for rows 6-17 the
# row-rotation caused by GateToken moving from the bottom of that block to the
# top, shifting the other coins down by one row).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '297.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.65%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.43'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.11%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.033'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.15%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07550'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.35%'

# Row 6
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.375'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.28%'

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.594'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.37%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9283'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.10%'

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.404'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.30%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1199'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.96%'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1842'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '7.02%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08886'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.23%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04038'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.04%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1054'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.02%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001286'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.97%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005801'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.29%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.339'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.60%'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3312'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.06%'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.988'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.39%'

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.08%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04054'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.95%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001264'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.37%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004165'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '7.08%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001229'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.98%'

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02421'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '4.04%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05210'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.97%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006527'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.89%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007795'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.16%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1330'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.28%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007551'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '2.46%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007830'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '10.85%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3224'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '11.54%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006784'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '5.77%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.09%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.004198'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-0.09%'

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.09%'

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.09%'
